$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert "Moving Time" column (C2:C12) from "h:mm:ss" text strings into
# plain numeric total-seconds values, since the API-derived time
# conversion is no longer needed now that Power BI handles formatting.

$movingTimes = @{
    2  = "01:34:05"
    3  = "00:46:54"
    4  = "01:35:30"
    5  = "00:36:20"
    6  = "01:21:54"
    7  = "01:16:24"
    8  = "00:20:36"
    9  = "01:12:11"
    10 = "00:47:31"
    11 = "00:44:31"
    12 = "01:44:35"
}

foreach ($row in $movingTimes.Keys) {
    $parts = $movingTimes[$row].Split(":")
    $hours = [int]$parts[0]
    $minutes = [int]$parts[1]
    $seconds = [int]$parts[2]
    $totalSeconds = $hours * 3600 + $minutes * 60 + $seconds

    $ws.Cells.Item($row, 3).Value = $totalSeconds
}
